$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '59.106.93'
$ws.Range('E2').Value = '  +1.20%  '

# Row 3
$ws.Range('D3').Value = '2.591.11'
$ws.Range('E3').Value = '  -0.01%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.19%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.54'
$ws.Range('E5').Value = '  +2.07%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.83'
$ws.Range('E6').Value = '  -1.33%  '

# Row 7
$ws.Range('E7').Value = '  +0.20%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  +0.22%  '

# Row 9
$ws.Range('D9').Value = '2.604.00'
$ws.Range('E9').Value = '  -0.37%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.44'
$ws.Range('E10').Value = '  -0.73%  '

# Row 11
$ws.Range('E11').Value = '  +0.60%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.333'
$ws.Range('E12').Value = '  +0.33%  '

# Row 13
$ws.Range('E13').Value = '  +2.85%  '

# Row 14
$ws.Range('D14').Value = '3.054.93'
$ws.Range('E14').Value = '  +0.10%  '

# Row 15
$ws.Range('D15').Value = '59.065.54'
$ws.Range('E15').Value = '  +1.14%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.46'
$ws.Range('E16').Value = '  +0.69%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.582.65'
$ws.Range('E17').Value = '  -0.41%  '

# Row 18
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000134'
$ws.Range('E18').Value = '  +0.18%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '347.35'
$ws.Range('E19').Value = '  +2.61%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.33'
$ws.Range('E20').Value = '  +0.34%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.10'
$ws.Range('E21').Value = '  -1.22%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.42'
$ws.Range('E22').Value = '  -0.08%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('E23').Value = '  -0.10%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.39'
$ws.Range('E24').Value = '  +3.01%  '

# Row 25
$ws.Range('E25').Value = '  -0.35%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.405'
$ws.Range('E26').Value = '  +0.57%  '

# Row 27
$ws.Range('E27').Value = '  +0.26%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.15'
$ws.Range('E28').Value = '  +1.55%  '

# Row 29
$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  +0.06%  '

# Row 30
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0736'
$ws.Range('E30').Value = '  -0.49%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.61'
$ws.Range('E31').Value = '  +2.58%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.86'
$ws.Range('E32').Value = '  -3.68%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.78'
$ws.Range('E33').Value = '  +0.05%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '148.87'
$ws.Range('E34').Value = '  -0.35%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.98'
$ws.Range('E35').Value = '  +0.02%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.12'
$ws.Range('E36').Value = '  -1.03%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.82'
$ws.Range('E37').Value = '  +1.69%  '

# Row 38
$ws.Range('E38').Value = '  +0.51%  '

# Row 39
$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.824'
$ws.Range('E39').Value = '  -1.01%  '

# Row 40
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.828'
$ws.Range('E40').Value = '  -3.35%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.53'
$ws.Range('E41').Value = '  +0.22%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'

# Row 43
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '268.97'
$ws.Range('E43').Value = '  -1.91%  '

# Row 44
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.595'
$ws.Range('E44').Value = '  -1.01%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.75'
$ws.Range('E45').Value = '  +0.62%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0958'
$ws.Range('E46').Value = '  +1.01%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0519'
$ws.Range('E47').Value = '  -0.32%  '

# Row 48
$ws.Range('E48').Value = '  -0.85%  '

# Row 49
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.958.95'
$ws.Range('E49').Value = '  -0.64%  '

# Row 50
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.62'
$ws.Range('E50').Value = '  -0.51%  '

# Row 51
$ws.Range('E51').Value = '  +0.41%  '
